$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert "Brasil / 01/01/2022" row before the Nordeste block (new row 14) ---
$ws.Rows.Item(14).Insert()
$ws.Range("A14").Value = "Brasil"
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "01/01/2022"
$ws.Range("B14").ClearFormats()
$ws.Range("C14").Value = 0.437868572310842

# --- Insert "Nordeste / 01/01/2022" row after the Nordeste block (new row 27) ---
$ws.Rows.Item(27).Insert()
$ws.Range("A27").Value = "Nordeste"
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "01/01/2022"
$ws.Range("B27").ClearFormats()
$ws.Range("C27").Value = 0.4255447605153426

# --- Append "Sergipe / 01/01/2022" row at the end (new row 40) ---
$ws.Range("A40").Value = "Sergipe"
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "01/01/2022"
$ws.Range("B40").ClearFormats()
$ws.Range("C40").Value = 0.3925077107324511
